$d = $word.ActiveDocument

# 1. Update the table caption number from "Supplemental Table 3." to "Supplemental Table 4."
$d.Content.Find.Execute("Supplemental Table 3.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Supplemental Table 4.", 2) | Out-Null

# 2. Append sentence about temperature variability definition to the caption text
$d.Content.Find.Execute("Associations of nestling mass and temperature, assessed in separate models stratified by relative nestling size at mid development measure (smallest vs. other).",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "Associations of nestling mass and temperature, assessed in separate models stratified by relative nestling size at mid development measure (smallest vs. other). Temperature variability is defined as the interquartile range.",
                         2) | Out-Null

# 3. Update table header text "Effect of temperature IQR" -> "Effect of temperature variability"
$d.Content.Find.Execute("Effect of temperature IQR", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Effect of temperature variability", 2) | Out-Null

# 4. Update footnote text for minimum temperature models: R^2^ -> R-squared
$d.Content.Find.Execute("R^2^ for adjusted minimum temperature models. Small size model: Marginal R^2^ = 0.34, Conditional R^2^ = 0.91; Other size model: Marginal R^2^ = 0.33, Conditional R^2^ = 0.85",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "R-squared for adjusted minimum temperature models. Small size model: Marginal R-squared = 0.34, Conditional R-squared = 0.91; Other size model: Marginal R-squared = 0.33, Conditional R-squared = 0.85",
                         2) | Out-Null

# 5. Update footnote text for maximum temperature models: R^2^ -> R-squared
$d.Content.Find.Execute("R^2^ for adjusted maximum temperature models. Small size model: Marginal R^2^ = 0.37, Conditional R^2^ = 0.92; Other size model: Marginal R^2^ = 0.32, Conditional R^2^ = 0.84",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "R-squared for adjusted maximum temperature models. Small size model: Marginal R-squared = 0.37, Conditional R-squared = 0.92; Other size model: Marginal R-squared = 0.32, Conditional R-squared = 0.84",
                         2) | Out-Null

# 6. Update footnote text for temperature IQR models: R^2^ -> R-squared and IQR -> variability
$d.Content.Find.Execute("R^2^ for adjusted temperature IQR models. Small size model: Marginal R^2^ = 0.47, Conditional R^2^ = 0.91; Other size model: Marginal R^2^ = 0.49, Conditional R^2^ = 0.84",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "R-squared for adjusted temperature variability models. Small size model: Marginal R-squared = 0.47, Conditional R-squared = 0.91; Other size model: Marginal R-squared = 0.49, Conditional R-squared = 0.84",
                         2) | Out-Null
